# ncp-gop-transect-info.xlsx edit
# commit: "titles, defns, d180, fixed catvars, zip inputs"
#
# On the Keywords sheet:
#  - replace the old "dissolved oxygen" keyword (row 7) with
#    "gross primary production"
#  - add a new keyword row (row 8): "Northeast U.S. Continental Shelf"
#    with thesaurus "NOAA Large Marine Ecosystems"
#  - leave the Keywords sheet selected/active (A7:B8), matching the
#    final editing focus of the workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Keywords")

$ws.Range("A8").Value = "Northeast U.S. Continental Shelf"
$ws.Range("B8").Value = "NOAA Large Marine Ecosystems"
$ws.Range("A7").Value = "gross primary production"

$ws.Activate()
$ws.Range("A7:B8").Select()
